# logboek fotosjaak.xlsx - apply "MySqlDatabaseClass.php aangepast en logboek
# bijgewerkt ook link test db gemaakt" edit.
#
# Summary of the change (week48 sheet):
#   - G7 gets a real "(D7-C7)" duration formula instead of being blank.
#   - C8/D8 switch from the HH:MM:SS style to an HH:MM style, and G8 gets a
#     "(D8-C8)" duration formula too.
#   - Two brand new activity rows are added (row 9: "vrijdag" 29 nov 2013 /
#     db connectie code..., row 10: continuation / De class
#     MySqlDatabaseClass afgemaakt...) plus a third new row (row 11: nieuwe
#     link gemaakt test db clas...) that has no duration cell.
#   - The "totaal" row moves from row 11 down to row 13 (row 12 stays empty)
#     and its SUM now covers G7:G10.
#   - The "totaal" worksheet picks the new week48!G13 total up automatically
#     through its existing formula.

$wb = $excel.ActiveWorkbook
$ws4 = $wb.Worksheets.Item("week48")
$ws5 = $wb.Worksheets.Item("totaal")

# --- make room for the two extra activity rows -----------------------------
# Old row 11 ("totaal" / SUM) needs to end up on row 13, with row 12 left
# empty, so insert two rows above it. New (still blank) rows 11 & 12 appear;
# row 12 is cleared below so it disappears from the sheet entirely again.
$ws4.Rows.Item(11).Resize(2).Insert()
$ws4.Rows.Item(12).Clear()

# --- row 7: give the duration column a real formula ------------------------
$ws4.Range("G7").Formula = "=(D7-C7)"

# --- row 8: begin/eind switch to HH:MM, duration formula added -------------
$ws4.Range("C8:D8").NumberFormat = "HH:MM"
$ws4.Range("G8").NumberFormat = "HH:MM:SS"
$ws4.Range("G8").Formula = "=(D8-C8)"

# --- row 9: new activity (vrijdag 29 nov 2013) ------------------------------
$ws4.Range("A9").Value = "vrijdag"
$ws4.Range("B9").Value = 41607
$ws4.Range("B9").NumberFormat = "DD/MM/YY"
$ws4.Range("C9").Value = 0.401388888888889
$ws4.Range("D9").Value = 0.416666666666667
$ws4.Range("C9:D9").NumberFormat = "HH:MM"
$ws4.Range("E9").Value = 3
$ws4.Range("F9").Value = "db connectie code gemaakt in de class MySqlDatabaseClass"
$ws4.Range("F9").NumberFormat = "@"
$ws4.Range("F9").HorizontalAlignment = -4130
$ws4.Range("G9").NumberFormat = "HH:MM:SS"
$ws4.Range("G9").Formula = "=(D9-C9)"
$ws4.Rows.Item(9).RowHeight = 25.35

# --- row 10: continuation activity ------------------------------------------
$ws4.Range("C10").Value = 0.417361111111111
$ws4.Range("D10").Value = 0.427083333333333
$ws4.Range("C10:D10").NumberFormat = "HH:MM"
$ws4.Range("E10").Value = 4
$ws4.Range("F10").Value = "De class MySqlDatabaseClass afgemaakt. Link naar gemaakt in link.php"
$ws4.Range("F10").HorizontalAlignment = -4130
$ws4.Range("G10").NumberFormat = "HH:MM:SS"
$ws4.Range("G10").Formula = "=(D10-C10)"
$ws4.Rows.Item(10).RowHeight = 25.35

# --- row 11: new activity, no duration cell ---------------------------------
$ws4.Range("C11").Value = 0.4375
$ws4.Range("D11").Value = 0.458333333333333
$ws4.Range("C11:D11").NumberFormat = "HH:MM:SS"
$ws4.Range("E11").Value = 5
$ws4.Range("F11").Value = "nieuwe link gemaakt test db clas. Op deze pagina een object gemaakt van de MySqlDatabaseClass."
$ws4.Range("F11").HorizontalAlignment = -4130
$ws4.Rows.Item(11).RowHeight = 37.3
$ws4.Range("G11").Clear()

# --- row 13: "totaal" row, now summing G7:G10 -------------------------------
$ws4.Range("G13").NumberFormat = "HH:MM:SS"
$ws4.Range("G13").Formula = "=SUM(G7:G10)"

# week48 stays the active sheet/tab, now with the selection parked on F15
$ws4.Activate()
$ws4.Range("F15").Select()

# the "totaal" sheet's week48 total (and the grand total) recompute on their
# own via the existing =week48!G13 / =SUM(B7:B8) formulas.
